# April 14th - include approximate mapping from reported to announced date
# Appends 5 new daily rows (2020-04-10 .. 2020-04-14) to the bottom of the
# existing data table on Sheet1, extending it from row 43 to row 48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date serial, confirmed, totalConfirmed, probable, totalProbable, cases,
# totalCases, recovered, totalRecovered, inHospitalNow, (K=totalBeenInHospital
# left blank), inIcu, deaths, totalDeaths, overseas, contact, investigating,
# community, established, tag
$newRows = @(
    @{ Row = 44; A = 43931; B = 23; C = 1015; D = 21; E = 268; F = 44; G = 1283; H = 56; I = 373; J = 16; L = 4; M = 1; N = 2; O = 513; P = 565; Q = 180; R = 26; S = 1283; T = "Manual" },
    @{ Row = 45; A = 43932; B = 20; C = 1035; D = 9;  E = 277; F = 29; G = 1312; H = 49; I = 422; J = 15; L = 5; M = 2; N = 4; O = 525; P = 604; Q = 144; R = 26; S = 1312; T = "Manual" },
    @{ Row = 46; A = 43933; B = 14; C = 1049; D = 4;  E = 281; F = 18; G = 1330; H = 49; I = 471; J = 14; L = 5; M = 0; N = 4; O = 532; P = 625; Q = 146; R = 27; S = 1330; T = "Manual" },
    @{ Row = 47; A = 43934; B = 15; C = 1064; D = 4;  E = 285; F = 19; G = 1349; H = 75; I = 546; J = 15; L = 4; M = 1; N = 5; O = 540; P = 634; Q = 148; R = 27; S = 1349; T = "Manual" },
    @{ Row = 48; A = 43935; B = 8;  C = 1072; D = 9;  E = 294; F = 17; G = 1366; H = 82; I = 628; J = 15; L = 3; M = 4; N = 9; O = 533; P = 656; Q = 150; R = 27; S = 1366; T = "Manual" }
)

$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A carries the date serial number, formatted the same way as the
    # rest of the table (style index 2 -> "yyyy-mm-dd HH:mm:ss UTC").
    $dateCell = $ws.Cells.Item($rowNum, 1)
    $dateCell.Value = $r.A
    $dateCell.NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"

    foreach ($col in $cols) {
        if ($col -eq "A") { continue }
        $cell = $ws.Range($col + $rowNum)
        $cell.Value = $r[$col]
    }
}
